$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added at the top of the historical table
# (row 39), pushing every existing record down by one row and appending
# the previously-last record (old row 87) as the new row 88.
$ws.Rows("39:39").Insert()

$ws.Range("A39").Value = 8
$ws.Range("B39").Value = "Terminal La Palmera de La Serena"
$ws.Range("C39").Value = "Coquimbo"
$ws.Range("D39").Value = 44483
$ws.Range("E39").Value = 4
$ws.Range("F39").Value = 100112044
$ws.Range("G39").Value = "Perejil"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 3300
$ws.Range("K39").Value = 1500
$ws.Range("L39").Value = 2000
$ws.Range("M39").Value = 1750
$ws.Range("N39").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O39").Value = "Provincia del Elquí"
$ws.Range("P39").Value = 1167
$ws.Range("Q39").Value = 1.5
$ws.Range("R39").Value = "Hortaliza"
